$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column H (Industries) to 0 for rows 31 through 94
$ws.Range("H31:H94").Value = 0
